$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The header in column U ("p24_3") is removed entirely - delete the whole column.
$ws.Range("U1").EntireColumn.Delete()

# Column T (previously empty except T52=1.5) now gets a value of 0 for every
# data row (2-83); row 52's old 1.5 is overwritten with 0 as well.
$ws.Range("T2:T83").Value = 0
